$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed cells per diff. Cells whose new text would otherwise be
# auto-parsed by Excel as a number are pinned to Text format first so the
# literal string (matching the inlineStr content in the source) is kept.

$ws.Range('D2').Value = '54.314.13'
$ws.Range('E2').Value = '  -7.92%  '
$ws.Range('D3').Value = '2.889.13'
$ws.Range('E3').Value = '  -10.70%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '473.05'
$ws.Range('E5').Value = '  -12.46%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '126.11'
$ws.Range('E6').Value = '  -7.83%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = '2.886.77'
$ws.Range('E8').Value = '  -10.75%  '
$ws.Range('E9').Value = '  -12.31%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.62'
$ws.Range('E10').Value = '  -13.43%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0973'
$ws.Range('E11').Value = '  -15.43%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.333'
$ws.Range('E12').Value = '  -15.67%  '
$ws.Range('E13').Value = '  -3.71%  '
$ws.Range('D14').Value = '3.384.41'
$ws.Range('E14').Value = '  -10.68%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '22.90'
$ws.Range('E15').Value = '  -11.94%  '
$ws.Range('D16').Value = '54.376.84'
$ws.Range('E16').Value = '  -7.84%  '
$ws.Range('D17').Value = '2.894.46'
$ws.Range('E17').Value = '  -10.48%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0000135'
$ws.Range('E18').Value = '  -15.51%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.53'
$ws.Range('E20').Value = '  -13.60%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.14'
$ws.Range('E21').Value = '  -13.86%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '307.99'
$ws.Range('E22').Value = '  -15.03%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.998'
$ws.Range('E23').Value = '  -0.21%  '
$ws.Range('E24').Value = '  -13.93%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '59.46'
$ws.Range('E25').Value = '  -15.71%  '
$ws.Range('E26').Value = '  +0.14%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.153'
$ws.Range('E27').Value = '  -9.98%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.998'
$ws.Range('E28').Value = '  -0.11%  '
$ws.Range('D29').Value = '0.0₃0828'
$ws.Range('E29').Value = '  -15.16%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.21'
$ws.Range('E30').Value = '  -12.53%  '
$ws.Range('E31').Value = '  -6.60%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.14'
$ws.Range('E32').Value = '  -13.76%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.15'
$ws.Range('E33').Value = '  -12.84%  '
$ws.Range('E34').Value = '  -16.57%  '
$ws.Range('B35').Value = 'NEARProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.25'
$ws.Range('E35').Value = '  -14.50%  '
$ws.Range('B36').Value = 'Monero'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '139.88'
$ws.Range('E36').Value = '  -14.28%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.46'
$ws.Range('E37').Value = '  -15.18%  '
$ws.Range('E38').Value = '  -15.60%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '23.01'
$ws.Range('E39').Value = '  -13.31%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0622'
$ws.Range('E40').Value = '  -12.76%  '
$ws.Range('D41').Value = '2.914.98'
$ws.Range('E41').Value = '  -10.79%  '
$ws.Range('E42').Value = '  +0.05%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '34.82'
$ws.Range('E43').Value = '  -15.54%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.968'
$ws.Range('E44').Value = '  -12.99%  '
$ws.Range('B45').Value = 'Filecoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.44'
$ws.Range('E45').Value = '  -14.85%  '
$ws.Range('B46').Value = 'Mantle'
$ws.Range('C46').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.597'
$ws.Range('E46').Value = '  -16.84%  '
$ws.Range('B47').Value = 'Stacks'
$ws.Range('C47').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.30'
$ws.Range('E47').Value = '  -13.61%  '
$ws.Range('B48').Value = 'Maker'
$ws.Range('C48').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D48').Value = '2.063.69'
$ws.Range('E48').Value = '  -10.47%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '18.16'
$ws.Range('E49').Value = '  -13.56%  '
$ws.Range('E50').Value = '  -15.65%  '
$ws.Range('E51').Value = '  -12.43%  '
